$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "tc024" right after "tc020" (before "tc013"),
#    mirroring the layout/format used on the "tc020" sheet, and add the
#    TC024 test data.
# ------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("tc020")
$newSheet = $wb.Worksheets.Add($null, $srcSheet)
$newSheet.Name = "tc024"

# Copy the formatting (font/style) used by tc020's header/value cells so the
# new sheet's cells share the same cell style as the other TC sheets.
$srcSheet.Range("A1:B2").Copy()
$newSheet.Range("A1:B2").PasteSpecial(-4122)

# Populate the new sheet's data (order matters so new shared strings land
# on the expected indexes).
$newSheet.Range("A1").Value = "Project Name"
$newSheet.Range("B1").Value = "release"
$newSheet.Range("A2").Value = "STG- PulseCodeOnAzureCloud"
$newSheet.Range("B2").Value = " Release update notoification 09-01-2026"
$newSheet.Range("C1").Value = "Status"
$newSheet.Range("C2").Value = "Planned"

# Column widths matching the other TC data sheets.
$newSheet.Columns.Item(1).ColumnWidth = 29.6165
$newSheet.Columns.Item(2).ColumnWidth = 33.33325

# Active cell / selection on the new sheet.
$newSheet.Range("C2").Select()

# ------------------------------------------------------------------
# 2. Update the selection on "tc020" (it used to point at A4, now the
#    whole A1:B2 block is selected).
# ------------------------------------------------------------------
$tc020 = $wb.Worksheets.Item("tc020")
$tc020.Range("A1:B2").Select()

# ------------------------------------------------------------------
# 3. Make the newly added "tc024" tab the active/selected tab.
# ------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("C2").Select()
